$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.043.57'
$ws.Range("E2").Value = '  -1.18%  '

$ws.Range("D3").Value = '2.636.21'
$ws.Range("E3").Value = '  +1.07%  '

$ws.Range("E4").Value = '  +0.09%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '606.12'
$ws.Range("E5").Value = '  +2.17%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.49'
$ws.Range("E6").Value = '  -1.35%  '

$ws.Range("E7").Value = '  +0.00%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.584'
$ws.Range("E8").Value = '  -0.91%  '

$ws.Range("D9").Value = '2.634.92'
$ws.Range("E9").Value = '  +1.13%  '

$ws.Range("E10").Value = '  +0.47%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.61'
$ws.Range("E11").Value = '  -0.67%  '

$ws.Range("E12").Value = '  +0.36%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.363'
$ws.Range("E13").Value = '  +2.70%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.08'
$ws.Range("E14").Value = '  -0.60%  '

$ws.Range("D15").Value = '3.110.61'
$ws.Range("E15").Value = '  +1.38%  '

$ws.Range("D16").Value = '62.882.02'
$ws.Range("E16").Value = '  -1.09%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000144'
$ws.Range("E17").Value = '  -1.39%  '

$ws.Range("D18").Value = '2.638.45'
$ws.Range("E18").Value = '  +2.69%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.33'
$ws.Range("E19").Value = '  +1.20%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.44'
$ws.Range("E20").Value = '  +2.11%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '339.70'
$ws.Range("E21").Value = '  -0.69%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.81'
$ws.Range("E22").Value = '  +1.11%  '

$ws.Range("E23").Value = '  -0.25%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '66.58'
$ws.Range("E24").Value = '  -2.89%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.63'
$ws.Range("E25").Value = '  +0.83%  '

$ws.Range("E26").Value = '  -2.38%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.68'
$ws.Range("E27").Value = '  +3.99%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.163'
$ws.Range("E28").Value = '  -1.37%  '

$ws.Range("B29").Value = 'Binance-PegBSC-USD'
$ws.Range("C29").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  +0.15%  '

$ws.Range("B30").Value = 'Bittensor'
$ws.Range("C30").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '532.35'
$ws.Range("E30").Value = '  +10.60%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.80'
$ws.Range("E31").Value = '  -1.69%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.03'
$ws.Range("E32").Value = '  +2.34%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.78'
$ws.Range("E33").Value = '  +2.95%  '

$ws.Range("D34").Value = '0.0₃0803'
$ws.Range("E34").Value = '  -1.82%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '170.90'
$ws.Range("E35").Value = '  -3.26%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.14'
$ws.Range("E36").Value = '  +13.44%  '

$ws.Range("B37").Value = 'FirstDigitalUSD'
$ws.Range("C37").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("E37").Value = '  -0.09%  '

$ws.Range("B38").Value = 'PolygonEcosystemToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.403'
$ws.Range("E38").Value = '  +0.95%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '19.05'
$ws.Range("E39").Value = '  +0.54%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.88'
$ws.Range("E40").Value = '  +9.03%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '172.63'
$ws.Range("E41").Value = '  +5.13%  '

$ws.Range("E42").Value = '  +0.09%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.72'
$ws.Range("E43").Value = '  -0.94%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '22.30'
$ws.Range("E44").Value = '  +2.40%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0569'
$ws.Range("E45").Value = '  +5.32%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.627'
$ws.Range("E46").Value = '  -0.23%  '

$ws.Range("B47").Value = 'VeChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0240'
$ws.Range("E47").Value = '  +0.45%  '

$ws.Range("B48").Value = 'Stellar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0960'
$ws.Range("E48").Value = '  -0.36%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '18.53'
$ws.Range("E49").Value = '  +0.22%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.73'
$ws.Range("E50").Value = '  -0.65%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '11.22'
$ws.Range("E51").Value = '  -1.31%  '
